$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the existing segment names from column A (rows 2-20) before shifting
$segments = @()
for ($r = 2; $r -le 20; $r++) {
    $segments += $ws.Cells.Item($r, 1).Value2
}

# Insert a new column before column B; this shifts old B:F -> C:G
$ws.Columns("B:B").Insert()

# Copy the header cell formatting (bold/centered/bordered style) from the
# neighboring header cell into the new header cell, then set its text.
$ws.Cells.Item(1, 3).Copy() | Out-Null
$ws.Cells.Item(1, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(1, 2).Value = "segments"

# Fill in segment names (text) into column B, and numeric index (0-based) into column A
for ($r = 2; $r -le 20; $r++) {
    $ws.Cells.Item($r, 2).Value = $segments[$r - 2]
    $ws.Cells.Item($r, 1).Value = $r - 2
}

# The column insert copied column A's style (bold/centered/bordered) into the
# new column B data cells; clear that back to the default "Normal" style so
# the segment-name cells are unstyled, matching the rest of the data columns.
$ws.Range("B2:B20").Style = "Normal"

$excel.CutCopyMode = 0
